# Updated cryptos list on Fri Jun  2 17:52:51 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "27.076.52"
$ws.Range("E2").Value = "  +0.35%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.893.49"
$ws.Range("E3").Value = "  +1.03%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.07"
$ws.Range("E5").Value = "  +0.49%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  +0.00%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5185"
$ws.Range("E7").Value = "  +2.44%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3759"
$ws.Range("E8").Value = "  +2.67%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07218"
$ws.Range("E9").Value = "  +0.21%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  +1.92%  "

# Row 11 - Polygon
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8975"
$ws.Range("E11").Value = "  +0.31%  "

# Row 12 - TRON
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07660"
$ws.Range("E12").Value = "  +1.74%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.891.97"
$ws.Range("E13").Value = "  +1.74%  "

# Row 14 - Litecoin
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "94.23"
$ws.Range("E14").Value = "  -0.50%  "

# Row 15 - Polkadot
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.230"
$ws.Range("E15").Value = "  -0.21%  "

# Row 16 - BinanceUSD
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.000"

# Row 17 - ShibaInu
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008519"
$ws.Range("E17").Value = "  -0.35%  "

# Row 18 - Avalanche
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.44"
$ws.Range("E18").Value = "  +1.36%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.03%  "

# Row 20 - WrappedBTC
$ws.Range("D20").Value = "27.128.79"
$ws.Range("E20").Value = "  +0.39%  "

# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.062"
$ws.Range("E21").Value = "  +0.61%  "

# Row 22 - WrappedliquidstakedEther2.0
$ws.Range("D22").Value = "2.120.14"
$ws.Range("E22").Value = "  +1.68%  "

# Row 23 - Cosmos
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.59"
$ws.Range("E23").Value = "  +1.77%  "

# Row 24 - Chainlink
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.412"
$ws.Range("E24").Value = "  -0.17%  "

# Row 25 - LidoDAOToken
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.283"
$ws.Range("E25").Value = "  +9.83%  "

# Row 26 - Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "146.61"
$ws.Range("E26").Value = "  -1.23%  "

# Row 27 - Toncoin
$ws.Range("E27").Value = "  -2.83%  "

# Row 28 - EthereumClassic
$ws.Range("E28").Value = "  +0.73%  "

# Row 29 - BitcoinCash
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.46"

# Row 30 - Filecoin
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.978"
$ws.Range("E30").Value = "  +6.11%  "

# Row 31 - InternetComputer(DFINITY)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.785"
$ws.Range("E31").Value = "  +1.65%  "

# Row 33 - Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05046"
$ws.Range("E33").Value = "  -1.85%  "

# Row 34 - ARBITRUM
$ws.Range("E34").Value = "  +6.57%  "

# Row 35 - ImmutableX
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7731"
$ws.Range("E35").Value = "  +2.80%  "

# Row 36 - HuobiToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.980"
$ws.Range("E36").Value = "  -0.24%  "

# Row 37 - MXToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.284"
$ws.Range("E37").Value = "  +1.79%  "

# Row 38 - RenderToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.594"
$ws.Range("E38").Value = "  +1.10%  "

# Row 39 - TheSandbox
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5608"
$ws.Range("E39").Value = "  -0.86%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -0.67%  "

# Row 41 - TrustWalletToken
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.075"
$ws.Range("E41").Value = "  +0.21%  "

# Row 42 - Aptos->Quant
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "119.51"
$ws.Range("E42").Value = "  +3.37%  "

# Row 43 - FraxShare
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "6.622"
$ws.Range("E43").Value = "  +0.14%  "

# Row 44 - Quant->Aptos
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.956"
$ws.Range("E44").Value = "  +4.98%  "

# Row 45 - Algorand
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1517"
$ws.Range("E45").Value = "  +2.80%  "

# Row 46 - Decentraland
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4820"
$ws.Range("E46").Value = "  +1.74%  "

# Row 47 - PaxDollar
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  +0.02%  "

# Row 48 - EnergySwap
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.14"
$ws.Range("E48").Value = "  -0.18%  "

# Row 49 - NEARProtocol
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.595"
$ws.Range("E49").Value = "  +1.87%  "

# Row 50 - Elrond
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.42"
$ws.Range("E50").Value = "  +1.40%  "

# Row 51 - Aave
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "63.86"
$ws.Range("E51").Value = "  +1.16%  "
